$d = $word.ActiveDocument

# 1) Replace the standalone [[PERSON_64]] placeholder (job position field) with literal text
$d.Content.Find.Execute("[[PERSON_64]]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Backend Developer", 2)

# 2) Renumber [[PERSON_65]] -> [[PERSON_64]] (keeps " Studio s.r.o." suffix intact)
$d.Content.Find.Execute("[[PERSON_65]] Studio s.r.o.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[[PERSON_64]] Studio s.r.o.", 2)

# 3) Renumber [[PERSON_66]] -> [[PERSON_65]] (keeps " & Partners s.r.o." suffix intact)
$d.Content.Find.Execute("[[PERSON_66]] & Partners s.r.o.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[[PERSON_65]] & Partners s.r.o.", 2)

# 4) Renumber [[PERSON_67]] -> [[PERSON_66]] (keeps "Kontaktní osoba: " prefix intact)
$d.Content.Find.Execute("Kontaktní osoba: [[PERSON_67]]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Kontaktní osoba: [[PERSON_66]]", 2)
